# Update the cryptocurrency price/volume snapshot to the latest scraped values.
# (GitHub Actions data refresh - see commit message)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.499.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "'3.816.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'700.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.84%  "

$ws.Range("D6").Value = "'174.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.39%  "

$ws.Range("D7").Value = "'3.811.99"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("D11").Value = "'7.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.69%  "

$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.39%  "

$ws.Range("D14").Value = "'36.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").Value = "'4.463.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").Value = "'3.816.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").Value = "'71.474.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").Value = "'17.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").Value = "'11.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.00%  "

$ws.Range("D22").Value = "'486.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("D23").Value = "'0.717"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").Value = "'84.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").Value = "'12.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("D29").Value = "'3.967.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'3.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.14%  "

$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("D33").Value = "'7.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.98%  "

$ws.Range("D34").Value = "'0.184"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "

$ws.Range("D35").Value = "'29.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").Value = "'9.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("D39").Value = "'2.40"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").Value = "'3.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").Value = "'6.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.53%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.70%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'163.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "

$ws.Range("D46").Value = "'0.000310"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.72%  "

$ws.Range("D47").Value = "'44.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("D48").Value = "'48.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("D49").Value = "'417.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.77%  "

$ws.Range("D50").Value = "'0.302"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").Value = "'8.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
